$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add "Sheet3" right after "Sheet2" and make it the active sheet.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Sheet3"

# ---------------------------------------------------------------------
# Student names (column A) - typed first, top to bottom.
# ---------------------------------------------------------------------
$ws3.Range("A2").Value = " Abhi"
$ws3.Range("A3").Value = " komal"
$ws3.Range("A4").Value = " Shivani"
$ws3.Range("A5").Value = " Taksh"
$ws3.Range("A6").Value = " Ixit"

# Header row for the marks table.
$ws3.Range("A1").Value = "NAME"
$ws3.Range("B1").Value = "BIO"
$ws3.Range("C1").Value = "MATHS"
$ws3.Range("D1").Value = "PHYSICS"

# ---------------------------------------------------------------------
# Raw marks (columns B, C, D) for rows 2-10.
# ---------------------------------------------------------------------
$bio = @(85, 96, 78, 56, 54, 30, 52, 63, 65)
$maths = @(56, 85, 75, 34, 85, 78, 98, 87, 75)
$physics = @(75, 89, 85, 68, 78, 78, 65, 85, 78)
for ($i = 0; $i -lt 9; $i++) {
    $r = 2 + $i
    $ws3.Range("B$r").Value = $bio[$i]
    $ws3.Range("C$r").Value = $maths[$i]
    $ws3.Range("D$r").Value = $physics[$i]
}

# ---------------------------------------------------------------------
# Statistic labels in column J (entered in the order the author typed
# them - note J7 precedes J6).
# ---------------------------------------------------------------------
$ws3.Range("J2").Value = "MORE THEN 80 IN MATHS"
$ws3.Range("J3").Value = "MORE THEN 80 IN BIO"
$ws3.Range("J4").Value = "MORE THEN 80 BIO & MAT."
$ws3.Range("J5").Value = "AVE. MARKES OF MATHS"
$ws3.Range("J7").Value = "AVE. MARKES OF BIO"
$ws3.Range("J6").Value = "AVE. MARKES OF PHYSICS"
$ws3.Range("J9").Value = "MORETHEN AVE. IN MATHS"
$ws3.Range("J10").Value = "MORE THEN AVE. IN BIO"
$ws3.Range("J11").Value = "MORE THEN AVE. IN PHYSICS"

# Remaining header-row labels.
$ws3.Range("E1").Value = "TOTAL"
$ws3.Range("F1").Value = "PER. OF STUDENT"
$ws3.Range("G1").Value = "MAX"

$ws3.Range("J13").Value = "HIGEST MARKES IN MATHS"
$ws3.Range("J14").Value = "HIGEST MARKES IN BIO"
$ws3.Range("J15").Value = "HIGEST MARKES IN PHYSICS"

$ws3.Range("C11").Value = "MAX "

$ws3.Range("H1").Value = "PASS OR FAIL"
$ws3.Range("I1").Value = "GREAD"

# ---------------------------------------------------------------------
# Per-student formulas (row 2 entered alone, rows 3:10 filled as one
# block so the engine records them as shared formulas, matching how
# Excel fills a formula down a selection).
# ---------------------------------------------------------------------
$ws3.Range("E2").Formula = "=SUM(B2:D2)"
$ws3.Range("E3:E10").Formula = "=SUM(B3:D3)"

$ws3.Range("F2:F10").NumberFormat = "0.00"
$ws3.Range("F2").Formula = "=E2/3"
$ws3.Range("F3:F10").Formula = "=E3/3"

$ws3.Range("G2").Formula = "=MAX(B2:D2)"
$ws3.Range("G3:G10").Formula = "=MAX(B3:D3)"

$ws3.Range("H2").Formula = '=IF(OR(B2<35,C2<35,D2<35),"FAIL","PASS")'
$ws3.Range("H3:H10").Formula = '=IF(OR(B3<35,C3<35,D3<35),"FAIL","PASS")'

$ws3.Range("I2").Formula = '=IF(H2="FAIL","FAIL",IF(F2>=90,"A",IF(F2>=60,"B",IF(F2>=50,"C"))))'
$ws3.Range("I3:I10").Formula = '=IF(H3="FAIL","FAIL",IF(F3>=90,"A",IF(F3>=60,"B",IF(F3>=50,"C"))))'

# ---------------------------------------------------------------------
# Summary stats in column K (count / average helpers).
# ---------------------------------------------------------------------
$ws3.Range("K2").NumberFormat = "0.00"
$ws3.Range("K2").Formula = '=COUNTIF(C2:C10,">80")'

$ws3.Range("K3").NumberFormat = "0.00"
$ws3.Range("K3").Formula = '=COUNTIF(B2:B10,">80")'

$ws3.Range("K4").NumberFormat = "0.00"
$ws3.Range("K4").Formula = '=COUNTIFS(C2:C10,">80",B2:B10,">80")'

$ws3.Range("K5").NumberFormat = "0.00"
$ws3.Range("K5").Formula = "=AVERAGE(C2:C10)"

$ws3.Range("K6").NumberFormat = "0.00"
$ws3.Range("K6").Formula = "=AVERAGE(D2:D10)"

$ws3.Range("K7").NumberFormat = "0.00"
$ws3.Range("K7").Formula = "=AVERAGE(B2:B10)"

$ws3.Range("K9").NumberFormat = "0.00"
$ws3.Range("K9").Formula = '=COUNTIF(C2:C10,">"&K5)'

$ws3.Range("K10").NumberFormat = "0.00"
$ws3.Range("K10").Formula = '=COUNTIF(B2:B10,">"&K7)'

$ws3.Range("K11").NumberFormat = "0.00"
$ws3.Range("K11").Formula = '=COUNTIF(D2:D10,">"&K6)'

# ---------------------------------------------------------------------
# Remaining names (rows 7-10) re-use names already used on Sheet2.
# ---------------------------------------------------------------------
$ws3.Range("A7").Value = "Himanshu"
$ws3.Range("A8").Value = "Meera"
$ws3.Range("A9").Value = "Dhara"
$ws3.Range("A10").Value = "Raj"

# ---------------------------------------------------------------------
# "MAX" row/column helpers.
# ---------------------------------------------------------------------
$ws3.Range("B11").Value = "MAX"
$ws3.Range("D11").Value = "MAX"

$ws3.Range("B12").Formula = "=MAX(B2:B10)"
$ws3.Range("C12").Formula = "=MAX(C2:C10)"
$ws3.Range("D12").Formula = "=MAX(D2:D10)"

$ws3.Range("F12").Value = "MAX"

$ws3.Range("F13").NumberFormat = "0.00"
$ws3.Range("F13").Formula = "=MAX(F2:F10)"

$ws3.Range("K13").NumberFormat = "0.00"
$ws3.Range("K13").Formula = "=LOOKUP(C12,C2:C10,A2:A10)"

$ws3.Range("B14").Formula = "=LOOKUP(B12,B2:B10,A2:A10)"

$ws3.Range("K14").NumberFormat = "0.00"
$ws3.Range("K14").Formula = "=LOOKUP(96,B2:B10,A2:A10)"

$ws3.Range("K15").NumberFormat = "0.00"
$ws3.Range("K15").Formula = "=LOOKUP(D12,D1:D10,A1:A10)"

# ---------------------------------------------------------------------
# Column widths (author widened/auto-fit these before finishing up).
# ---------------------------------------------------------------------
$ws3.Columns.Item(1).AutoFit() | Out-Null
$ws3.Columns.Item(6).AutoFit() | Out-Null
$ws3.Columns.Item(8).AutoFit() | Out-Null
$ws3.Columns.Item(9).ColumnWidth = 12.5703125
$ws3.Columns.Item(10).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Selections: Sheet2 loses the "active" flag and keeps a different
# selected cell; Sheet3 becomes the active/selected sheet.
# ---------------------------------------------------------------------
$ws2.Range("H13").Select()
$ws3.Activate()
$ws3.Range("F14").Select()
